$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 167.93333
$ws.Range("I42").Value = 165.08333
$ws.Range("J42").Value = 179.33333
$ws.Range("K42").Value = 495.24999
$ws.Range("L42").Value = 537.99999
$ws.Range("M42").Value = -265.24999
$ws.Range("N42").Value = -997.99999
$ws.Range("H69").Value = 7650.353
$ws.Range("I69").Value = 5912
$ws.Range("K69").Value = 17736
$ws.Range("M69").Value = -16862
$ws.Range("H72").Value = 7650.353
$ws.Range("I72").Value = 5912
$ws.Range("K72").Value = 53208
$ws.Range("M72").Value = -48840
$ws.Range("H74").Value = 11609
$ws.Range("I74").Value = 5566.6665
$ws.Range("K74").Value = 5566.6665
$ws.Range("M74").Value = -4630.6665
$ws.Range("H77").Value = 11609
$ws.Range("I77").Value = 5566.6665
$ws.Range("K77").Value = 27833.3325
$ws.Range("M77").Value = -23153.3325
$ws.Range("H104").Value = 218.75
$ws.Range("I104").Value = 94.333336
$ws.Range("K104").Value = 283.000008
$ws.Range("M104").Value = 1463.999992
$ws.Range("H113").Value = 7266.5
$ws.Range("I113").Value = 4379.857
$ws.Range("K113").Value = 4379.857
$ws.Range("M113").Value = -1125.857
$ws.Range("H137").Value = 2833.25
$ws.Range("I137").Value = 1764.8636
$ws.Range("J137").Value = 3392.8809
$ws.Range("K137").Value = 5294.5908
$ws.Range("L137").Value = 10178.6427
$ws.Range("M137").Value = -2744.5908
$ws.Range("N137").Value = -15278.6427

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5186.6787
$ws.Range("I61").Value = 4304.185
$ws.Range("K61").Value = 4304.185
$ws.Range("M61").Value = -4092.185
$ws.Range("H74").Value = 25643472
$ws.Range("I74").Value = 27780136
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 27780136
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -27779262
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 25643472
$ws.Range("I77").Value = 27780136
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 138900680
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -138896312
$ws.Range("N77").Value = -26236
$ws.Range("H110").Value = 3824.75
$ws.Range("I110").Value = 3181
$ws.Range("J110").Value = 8331
$ws.Range("K110").Value = 3181
$ws.Range("L110").Value = 8331
$ws.Range("M110").Value = -1136
$ws.Range("N110").Value = -12421
$ws.Range("H122").Value = 2743.3225
$ws.Range("I122").Value = 2233.72
$ws.Range("K122").Value = 6701.16
$ws.Range("M122").Value = -4251.16
$ws.Range("H132").Value = 2695.5518
$ws.Range("I132").Value = 1967.7693
$ws.Range("J132").Value = 9003
$ws.Range("K132").Value = 5903.3079
$ws.Range("L132").Value = 27009
$ws.Range("M132").Value = -3373.3079
$ws.Range("N132").Value = -32069
$ws.Range("H136").Value = 5186.6787
$ws.Range("I136").Value = 4304.185
$ws.Range("K136").Value = 12912.555
$ws.Range("M136").Value = -10362.555

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4621.048
$ws.Range("J20").Value = 6094.4
$ws.Range("L20").Value = 6094.4
$ws.Range("N20").Value = -6588.4
$ws.Range("H134").Value = 4881.6665
$ws.Range("I134").Value = 4325.077
$ws.Range("J134").Value = 8499.5
$ws.Range("K134").Value = 12975.231
$ws.Range("L134").Value = 25498.5
$ws.Range("M134").Value = -10440.231
$ws.Range("N134").Value = -30568.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 69186.164
$ws.Range("J135").Value = 69186.164
$ws.Range("L135").Value = 69186.164
$ws.Range("N135").Value = -79326.164

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 1999.5
$ws.Range("I76").Value = 1999.5
$ws.Range("K76").Value = 5998.5
$ws.Range("M76").Value = -5615.5
$ws.Range("H79").Value = 1999.5
$ws.Range("I79").Value = 1999.5
$ws.Range("K79").Value = 5998.5
$ws.Range("M79").Value = -4672.5
$ws.Range("H94").Value = 3160.1538
$ws.Range("I94").Value = 920
$ws.Range("K94").Value = 2760
$ws.Range("M94").Value = -2084
$ws.Range("H121").Value = 918
$ws.Range("I121").Value = 849.75
$ws.Range("K121").Value = 2549.25
$ws.Range("M121").Value = -1239.25
$ws.Range("H131").Value = 10306404
$ws.Range("I131").Value = 17858202
$ws.Range("J131").Value = 8103797
$ws.Range("K131").Value = 53574606
$ws.Range("L131").Value = 24311391
$ws.Range("M131").Value = -53569566
$ws.Range("N131").Value = -24321471
$ws.Range("H134").Value = 13435.143
$ws.Range("I134").Value = 18565.666
$ws.Range("K134").Value = 55696.99800000001
$ws.Range("M134").Value = -50626.99800000001
$ws.Range("H137").Value = 77562.08
$ws.Range("I137").Value = 692.1667
$ws.Range("J137").Value = 1000001
$ws.Range("K137").Value = 2076.5001
$ws.Range("L137").Value = 3000003
$ws.Range("M137").Value = 3023.4999
$ws.Range("N137").Value = -3010203

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 27499
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 27499
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 27499
$ws.Range("N44").Value = -28691
$ws.Range("M44").ClearContents()
$ws.Range("H107").Value = 608.35
$ws.Range("I107").Value = 677.4666999999999
$ws.Range("J107").Value = 401
$ws.Range("K107").Value = 677.4666999999999
$ws.Range("L107").Value = 401
$ws.Range("M107").Value = 1242.5333
$ws.Range("N107").Value = -4241
$ws.Range("H113").Value = 3096.3845
$ws.Range("I113").Value = 2601.6
$ws.Range("K113").Value = 2601.6
$ws.Range("M113").Value = -431.5999999999999
$ws.Range("H132").Value = 2401.652
$ws.Range("I132").Value = 1300.625
$ws.Range("J132").Value = 4918.2856
$ws.Range("K132").Value = 3901.875
$ws.Range("L132").Value = 14754.8568
$ws.Range("M132").Value = -1371.875
$ws.Range("N132").Value = -19814.8568

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3100.7666
$ws.Range("I22").Value = 2394.7058
$ws.Range("J22").Value = 4024.077
$ws.Range("K22").Value = 2394.7058
$ws.Range("L22").Value = 4024.077
$ws.Range("M22").Value = -2099.7058
$ws.Range("N22").Value = -4614.077
$ws.Range("H27").Value = 3100.7666
$ws.Range("I27").Value = 2394.7058
$ws.Range("J27").Value = 4024.077
$ws.Range("K27").Value = 2394.7058
$ws.Range("L27").Value = 4024.077
$ws.Range("M27").Value = -2287.7058
$ws.Range("N27").Value = -4238.077
$ws.Range("H55").Value = 1697.1724
$ws.Range("I55").Value = 591.1667
$ws.Range("J55").Value = 2477.8823
$ws.Range("K55").Value = 591.1667
$ws.Range("L55").Value = 2477.8823
$ws.Range("M55").Value = -418.1667
$ws.Range("N55").Value = -2823.8823
$ws.Range("H132").Value = 3872.9524
$ws.Range("I132").Value = 2956.24
$ws.Range("J132").Value = 5221.0586
$ws.Range("K132").Value = 8868.719999999999
$ws.Range("L132").Value = 15663.1758
$ws.Range("M132").Value = -6338.719999999999
$ws.Range("N132").Value = -20723.1758

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4979.6113
$ws.Range("I81").Value = 3244.3333
$ws.Range("J81").Value = 5326.6665
$ws.Range("K81").Value = 6488.6666
$ws.Range("L81").Value = 10653.333
$ws.Range("M81").Value = -5427.6666
$ws.Range("N81").Value = -12775.333
$ws.Range("H84").Value = 4979.6113
$ws.Range("I84").Value = 3244.3333
$ws.Range("J84").Value = 5326.6665
$ws.Range("K84").Value = 32443.333
$ws.Range("L84").Value = 53266.665
$ws.Range("M84").Value = -27139.333
$ws.Range("N84").Value = -63874.665
$ws.Range("H107").Value = 1495.381
$ws.Range("I107").Value = 1706.7059
$ws.Range("J107").Value = 597.25
$ws.Range("K107").Value = 5120.1177
$ws.Range("L107").Value = 1791.75
$ws.Range("M107").Value = -3200.1177
$ws.Range("N107").Value = -5631.75
$ws.Range("H113").Value = 1329.1111
$ws.Range("I113").Value = 1452.2727
$ws.Range("J113").Value = 1135.5714
$ws.Range("K113").Value = 4356.8181
$ws.Range("L113").Value = 3406.7142
$ws.Range("M113").Value = -2186.8181
$ws.Range("N113").Value = -7746.7142
$ws.Range("H132").Value = 1447.2222
$ws.Range("I132").Value = 1567.8572
$ws.Range("K132").Value = 4703.571599999999
$ws.Range("M132").Value = -2173.571599999999
